# Append newly-scraped departure rows (141-147) to the "Main Data" sheet,
# mirroring data pulled for an additional (11th) airport.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 141
$ws.Range("A141").Value = 140
$ws.Range("B141").Value = "Sunday, Jan 08"
$ws.Range("C141").Value = "8:20 PM"
$ws.Range("D141").Value = "UNKNOWN"
$ws.Range("E141").Value = "Budapest"
$ws.Range("F141").Value = "(BUD)"
$ws.Range("G141").Value = "Wizz Air "
$ws.Range("H141").Value = "A21N"
$ws.Range("I141").Value = "(HA-LZN)"
$ws.Range("J141").Value = "8:44 PM"
$ws.Range("L141").Value = "0 hours, 24 minutes"

# Row 142
$ws.Range("A142").Value = 141
$ws.Range("B142").Value = "Sunday, Jan 08"
$ws.Range("C142").Value = "8:35 PM"
$ws.Range("D142").Value = "FR7696"
$ws.Range("E142").Value = "Porto"
$ws.Range("F142").Value = "(OPO)"
$ws.Range("G142").Value = "Ryanair "
$ws.Range("H142").Value = "B738"
$ws.Range("I142").Value = "(EI-GJT)"
$ws.Range("J142").Value = "8:54 PM"
$ws.Range("L142").Value = "0 hours, 19 minutes"

# Row 143
$ws.Range("A143").Value = 142
$ws.Range("B143").Value = "Sunday, Jan 08"
$ws.Range("C143").Value = "8:35 PM"
$ws.Range("D143").Value = "U23818"
$ws.Range("E143").Value = "Paris"
$ws.Range("F143").Value = "(CDG)"
$ws.Range("G143").Value = "easyJet "
$ws.Range("H143").Value = "A320"
$ws.Range("I143").Value = "(OE-IVD)"
$ws.Range("J143").Value = "9:00 PM"
$ws.Range("L143").Value = "0 hours, 25 minutes"

# Row 144
$ws.Range("A144").Value = 143
$ws.Range("B144").Value = "Sunday, Jan 08"
$ws.Range("C144").Value = "8:55 PM"
$ws.Range("D144").Value = "FR3614"
$ws.Range("E144").Value = "Milan"
$ws.Range("F144").Value = "(BGY)"
$ws.Range("G144").Value = "Ryanair "
$ws.Range("H144").Value = "B738"
$ws.Range("I144").Value = "(SP-RSY)"
$ws.Range("J144").Value = "9:04 PM"
$ws.Range("L144").Value = "0 hours, 9 minutes"

# Row 145
$ws.Range("A145").Value = 144
$ws.Range("B145").Value = "Sunday, Jan 08"
$ws.Range("C145").Value = "9:25 PM"
$ws.Range("D145").Value = "LO3922"
$ws.Range("E145").Value = "Warsaw"
$ws.Range("F145").Value = "(WAW)"
$ws.Range("G145").Value = "LOT "
$ws.Range("H145").Value = "E170"
$ws.Range("I145").Value = "(SP-LDF)"
$ws.Range("J145").Value = "9:31 PM"
$ws.Range("L145").Value = "0 hours, 6 minutes"

# Row 146
$ws.Range("A146").Value = 145
$ws.Range("B146").Value = "Sunday, Jan 08"
$ws.Range("C146").Value = "9:45 PM"
$ws.Range("D146").Value = "FR6264"
$ws.Range("E146").Value = "Poznan"
$ws.Range("F146").Value = "(POZ)"
$ws.Range("G146").Value = "Ryanair "
$ws.Range("H146").Value = "B738"
$ws.Range("I146").Value = "(SP-RSM)"
$ws.Range("J146").Value = "9:47 PM"
$ws.Range("L146").Value = "0 hours, 2 minutes"

# Row 147
$ws.Range("A147").Value = 146
$ws.Range("B147").Value = "Sunday, Jan 08"
$ws.Range("C147").Value = "9:50 PM"
$ws.Range("D147").Value = "DY1043"
$ws.Range("E147").Value = "Oslo"
$ws.Range("F147").Value = "(OSL)"
$ws.Range("G147").Value = "Norwegian "
$ws.Range("H147").Value = "B738"
$ws.Range("I147").Value = "(LN-NIH)"
$ws.Range("J147").Value = "10:05 PM"
$ws.Range("L147").Value = "0 hours, 15 minutes"
